$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the policy table block (P3:T34) into new columns W:AA,
# replicating both the text content and the cell formatting (yellow
# highlight rows) of the existing policy table.
$src = $ws.Range("P3:T34")
$dst = $ws.Range("W3:AA34")
$src.Copy($dst)

# Row 8 gets an updated (v3) policy decision for this state: it now
# resolves to "change_lane" instead of "keep_distance" - reuse the
# existing shared string text from another cell that already reads
# "change_lane" for this same truth value.
$ws.Range("AA8").Value = $ws.Range("T16").Value2

# Update the current selection to reflect the newly added block.
$ws.Range("W3:AB34").Select()
